# faturamento_diario_lojas.xlsx - atualizacao dos dados bibi e add
# Rebuild rows 2-6 (store rows + total) with the updated figures:
#   row2 Bibi Cell Mundi        B=0       C=17144.8
#   row3 Bibi Cell Ponta Negra  B=1800.01 C=4670
#   row4 Bibi Cell Manauara     B=3340    C=1374
#   row5 Bibi Cell Vieiralves   B=0       C=4464
#   row6 total                  B=5140.01 C=27652.8
# All other daily columns (D..AF) stay 0, AG is the row total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; Name = "Bibi Cell Mundi";       B = 0;       C = 17144.8 },
    @{ Row = 3; Name = "Bibi Cell Ponta Negra"; B = 1800.01; C = 4670 },
    @{ Row = 4; Name = "Bibi Cell Manauara";    B = 3340;    C = 1374 },
    @{ Row = 5; Name = "Bibi Cell Vieiralves";  B = 0;       C = 4464 },
    @{ Row = 6; Name = "total";                 B = 5140.01; C = 27652.8 }
)

foreach ($r in $rows) {
    $rowIdx = $r.Row

    # Label in column A (bold/boxed style already applied on row2-4; extend it to the new rows too)
    $ws.Cells.Item($rowIdx, 1).Value = $r.Name
    if ($rowIdx -gt 4) {
        # Rows 5 and 6 are brand new - clone the existing label formatting (bold, border, centered)
        # from A2 instead of touching the "Normal" named style (which would wipe direct formatting).
        $ws.Cells.Item(2, 1).Copy()
        $ws.Cells.Item($rowIdx, 1).PasteSpecial(-4122)
    }

    # Column B and C carry the actual figures
    $ws.Cells.Item($rowIdx, 2).Value = $r.B
    $ws.Cells.Item($rowIdx, 3).Value = $r.C

    # Columns D..AF (4..32) are all zero
    for ($col = 4; $col -le 32; $col++) {
        $ws.Cells.Item($rowIdx, $col).Value = 0
    }

    # Column AG (33) = row total
    $ws.Cells.Item($rowIdx, 33).Value = $r.B + $r.C
}
